# Refresh crypto price/volume(1h) data as of the Fri Oct 25 17:13:50 UTC 2024 run.
# Numeric-looking Price (column D) values are written with a leading single-quote so
# Excel stores them as literal text (matching the original inlineStr cells) instead of
# auto-converting them to numbers; the leading quote itself is not stored in the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.759.43"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "2.533.52"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'590.84"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "'171.84"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("D9").Value = "2.531.46"
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E13").Value = "  -2.59%  "

$ws.Range("D14").Value = "'26.22"
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("D15").Value = "2.988.99"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("E16").Value = "  -1.67%  "

$ws.Range("D17").Value = "67.595.00"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").Value = "2.576.60"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").Value = "'8.09"
$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("E20").Value = "  +3.44%  "

$ws.Range("D21").Value = "'365.25"
$ws.Range("E21").Value = "  +1.03%  "

# Row 22/23: Binance-PegBSC-USD and Polkadot swapped positions, each with refreshed data.
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "'4.16"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("B23").Value = "Binance-PegBSC-USD"
$ws.Range("C23").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D23").Value = "'1.35"
$ws.Range("E23").Value = "  +35.48%  "

$ws.Range("D24").Value = "'4.52"
$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").Value = "'71.94"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  -6.05%  "

$ws.Range("D28").Value = "'9.77"
$ws.Range("E28").Value = "  -4.81%  "

$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("D31").Value = "'534.32"
$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("D32").Value = "'8.24"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("E34").Value = "  -3.62%  "

$ws.Range("E35").Value = "  -1.51%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'157.89"
$ws.Range("E37").Value = "  +1.40%  "

$ws.Range("D38").Value = "'19.41"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("E39").Value = "  -1.91%  "

$ws.Range("D40").Value = "'18.63"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "'1.77"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("D42").Value = "'5.06"
$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("E45").Value = "  -2.96%  "

$ws.Range("D46").Value = "'39.39"
$ws.Range("E46").Value = "  -1.10%  "

$ws.Range("D47").Value = "'148.48"
$ws.Range("E47").Value = "  +1.40%  "

$ws.Range("D48").Value = "'3.70"
$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").Value = "'0.549"
$ws.Range("E49").Value = "  -1.96%  "

$ws.Range("D50").Value = "0.0₆0272"
$ws.Range("E50").Value = "  -2.05%  "

$ws.Range("D51").Value = "'1.70"
$ws.Range("E51").Value = "  +1.08%  "
